$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F4").Value = "Test concatenate concatenateTest"
